$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Date"
$ws.Range("B2").Value = "Commodity/Service"
$ws.Range("C2").Value = "Group"
$ws.Range("D2").Value = "Price"
